# chore: adapt column header formatting to respective input file names
#
# The sheet's header row (row 1) used a generic "_old" / "_new" suffix to
# distinguish the two compared AHB format versions. Rename those headers so
# they carry the concrete format-version name instead ("_FV2410" / "_FV2504"),
# then turn the data range into a proper Excel Table and freeze the header
# row so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header labels: "*_old" -> "*_FV2410", "*_new" -> "*_FV2504" ---
$usedRange = $ws.UsedRange
$lastColumn = $usedRange.Columns.Count
$lastRow = $usedRange.Rows.Count

for ($col = 1; $col -le $lastColumn; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $headerText = $cell.Value2

    if ($headerText -ne $null) {
        if ($headerText.EndsWith("_old")) {
            $cell.Value = $headerText.Substring(0, $headerText.Length - 4) + "_FV2410"
        }
        elseif ($headerText.EndsWith("_new")) {
            $cell.Value = $headerText.Substring(0, $headerText.Length - 4) + "_FV2504"
        }
    }
}

# --- 2. Wrap the data range in a named Excel Table (ListObject) ---
$dataRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastColumn))
$listObject = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$listObject.Name = "Table1"
$listObject.TableStyle = ""

# --- 3. Freeze the header row (split below row 1) ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
